$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Packages")

# New PackageID rows appended to the Packages sheet (state_counters + packages refresh)
$newRows = @(
    @("SPA",    "MD", "Medicaid SPA", "",          "MD-25-9434",     "Under Review",    ""),
    @("SPA",    "MD", "Medicaid SPA", "",          "MD-25-9435",     "Approved",        ""),
    @("SPA",    "MD", "Medicaid SPA", "",          "MD-25-9436",     "Submitted",       ""),
    @("SPA",    "MD", "Medicaid SPA", "",          "MD-25-9437",     "Disapproved",     ""),
    @("Waiver", "MD", "1915(c)",      "Amendment", "MD-2260.R00.36", "Pending-Approval", "MD-2260.R00.00"),
    @("SPA",    "MD", "Medicaid SPA", "",          "MD-25-9438",     "Submitted",       "")
)

$startRow = 159
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $row = $newRows[$i]
    for ($c = 1; $c -le 7; $c++) {
        $cell = $ws.Cells.Item($r, $c)
        $text = $row[$c - 1]
        if ($text -eq "") {
            # A true blank value clears the cell outright under COM semantics,
            # so force an empty *text* entry (quote-prefixed empty string),
            # then strip the quote-prefix formatting back to the sheet default.
            $cell.Formula = "'"
            $cell.Style = "Normal"
        } else {
            $cell.Value = $text
        }
    }
}
